# Add "2022-Q3" data to the 港股/02202-万科企业 workbook:
#  1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#     existing "2022-Q2" sheet), populated with the fund-holding detail rows.
#  2. Insert a new summary row for "2022-Q3" at the top of the "总计" sheet's
#     data (row 2), pushing the existing quarters down by one row and
#     renumbering the index column (A) accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" worksheet, positioned before "2022-Q2"
# (item index 2 in the current, pre-insert, tab order).
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q3"

# Borrow header/format from the sheet that used to be "2022-Q2" (now shifted
# to item 3) so the new sheet's styling (bold/bordered header row, centered
# index column) matches the rest of the workbook exactly.
$template = $wb.Worksheets.Item(3)
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A19").PasteSpecial(-4122)

# Header labels
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Fund rows: 基金代码 / 基金名称 / 基金规模 / 股票总仓位 / 仓位占比 / 持有市值(亿元) / 仓位排名
$data = @(
  @("009100","安信稳健增利混合A","113.23","23.20","2.19","2.4797",1),
  @("009101","安信稳健增利混合C","84.02","23.20","2.19","1.8400",1),
  @("012609","安信稳健汇利一年持有混合A","35.76","23.15","2.33","0.8332",1),
  @("010864","泓德卓远混合A","22.84","92.87","2.92","0.6669",10),
  @("012256","安信丰穗一年持有混合A","24.55","24.64","2.16","0.5303",2),
  @("009849","安信稳健聚申一年持有期混合A","12.83","39.19","3.60","0.4619",2),
  @("012610","安信稳健汇利一年持有混合C","18.98","23.15","2.33","0.4422",1),
  @("008809","安信民稳增长混合A","9.52","49.26","4.29","0.4084",3),
  @("008810","安信民稳增长混合C","7.95","49.26","4.29","0.3411",3),
  @("010865","泓德卓远混合C","10.32","92.87","2.92","0.3013",10),
  @("671010","西部利得策略优选混合A","3.33","93.42","9.04","0.3010",3),
  @("012702","安信民安回报一年持有混合C","11.25","22.89","2.34","0.2632",1),
  @("012250","安信平衡增利混合A","2.63","60.12","5.91","0.1554",1),
  @("012251","安信平衡增利混合C","2.10","60.12","5.91","0.1241",1),
  @("012701","安信民安回报一年持有混合A","3.70","22.89","2.34","0.0866",1),
  @("010661","安信稳健聚申一年持有期混合C","2.28","39.19","3.60","0.0821",2),
  @("012257","安信丰穗一年持有混合C","2.43","24.64","2.16","0.0525",2),
  @("011060","西部利得策略优选混合C","0.47","93.42","9.04","0.0425",3)
)

# A blank, never-touched cell donates "no explicit style" formatting so the
# text-forced columns below don't end up tagged with a stray numeric-format
# style (mirrors the plain, un-styled B..G data cells used elsewhere in the
# workbook).
$blank = $newSheet.Range("Z100")

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $newSheet.Cells.Item($row,1).Value = $i

    # Columns B-G are stored as text in this workbook (even the numeric-
    # looking ones), so force text storage via NumberFormat "@" and then
    # strip the format change back off so no extra style is introduced.
    for ($col = 2; $col -le 7; $col++) {
        $cell = $newSheet.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rec[$col - 2]
        $blank.Copy()
        $cell.PasteSpecial(-4122)
    }

    # Column H (仓位排名) is a genuine number.
    $newSheet.Cells.Item($row,8).Value = $rec[6]
}

# ---------------------------------------------------------------------
# Step 2: insert the "2022-Q3" summary row at the top of "总计"'s table and
# shift the existing quarters down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Capture the existing rows 2-6 (in top-to-bottom order) before overwriting.
$existing = @(
  @("2022-Q2", 10, 1.84),
  @("2022-Q1", 7, 3.13),
  @("2021-Q4", 3, 1.77),
  @("2021-Q3", 4, 2),
  @("2021-Q1", 7, 3.07)
)

# Row 7 is brand new -- borrow the index column's style from row 6 first.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

for ($i = $existing.Length - 1; $i -ge 0; $i--) {
    $row = $i + 3
    $rec = $existing[$i]
    $total.Cells.Item($row,1).Value = $i + 1
    $total.Cells.Item($row,2).Value = $rec[0]
    $total.Cells.Item($row,3).Value = $rec[1]
    $total.Cells.Item($row,4).Value = $rec[2]
}

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 18
$total.Cells.Item(2,4).Value = 9.41

# Restore "总计" as the active sheet/tab (adding a worksheet makes the new
# sheet active by default).
$total.Activate()
